$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 39: Remove Nth Node from End of List -> fill in Status/Last Solved/Notes, update Topics ---
$ws.Range("H39").Value = "19/06/2025"
$ws.Range("D39").Value = "Linked Lists, Two Pointers"
$ws.Range("G39").Value = "STRUGGLED"
$ws.Range("I39").Value = "After looking at the solution it's not that bad. Clearly my mind still needs to get used to this data structure!"
$ws.Range("I39").WrapText = $true
$ws.Rows.Item(39).RowHeight = 30

# --- Row 40: Copy List with Random Pointer ---
$ws.Range("A40").Value = "Leetcode"
$ws.Range("B40").Value = 138
$ws.Range("B40").HorizontalAlignment = -4131
$ws.Range("C40").Value = "Copy List with Random Pointer"
$ws.Range("D40").Value = "Linked List"
$ws.Range("E40").Value = "Medium"
$ws.Range("F40").Value = "Neetcode 150"
$ws.Range("G40").Value = "STRUGGLED"
$ws.Range("H40").Value = "20/06/2025"

# --- Row 41: Add Two Numbers ---
$ws.Range("A41").Value = "Leetcode"
$ws.Range("B41").Value = 2
$ws.Range("B41").HorizontalAlignment = -4131
$ws.Range("C41").Value = "Add Two Numbers"
$ws.Range("D41").Value = "Linked List"
$ws.Range("E41").Value = "Medium"
$ws.Range("F41").Value = "Neetcode 150"
$ws.Range("G41").Value = "SOLVED"
$ws.Range("H41").Value = "20/06/2025"
$ws.Range("I41").Value = "I actually enjoyed this one!"
$ws.Range("I41").WrapText = $true

# --- View state: select the row after the new data (mirrors post-entry selection) ---
$ws.Range("A42:XFD42").Select()
